$d = $word.ActiveDocument

$pairs = @(
    @("424÷4=", "371÷2="),
    @("858÷5=", "734÷6="),
    @("527÷7=", "120÷6="),
    @("259÷9=", "353÷2="),
    @("729÷3=", "853÷7="),
    @("312÷9=", "334÷8="),
    @("582÷3=", "816÷8="),
    @("378÷4=", "312÷2="),
    @("430÷6=", "626÷2="),
    @("767÷3=", "803÷6="),
    @("413÷6=", "129÷7="),
    @("416÷4=", "264÷2="),
    @("376÷5=", "607÷8="),
    @("584÷5=", "636÷9="),
    @("688÷5=", "524÷6="),
    @("201÷6=", "352÷4="),
    @("279÷8=", "392÷5="),
    @("142÷8=", "992÷4="),
    @("750÷2=", "987÷3="),
    @("746÷8=", "674÷3="),
    @("539÷8=", "970÷6="),
    @("342÷3=", "231÷3="),
    @("609÷7=", "670÷6="),
    @("399÷6=", "640÷4="),
    @("513÷8=", "933÷6=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
